$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing data rows down by one.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New row 2 values
$ws.Range("A2").Value = -0.0303905457258224
$ws.Range("B2").Value = -0.02122756652534
$ws.Range("C2").Value = -0.038026362657547

# Append new rows 23-31 (row 22 already holds the former last row's data).
$newRows = @(
    @(-0.0070249503478407, 4.270253658294678, -0.0296269636601209),
    @(0.4137084782123565, 2.936276435852051, 0.2823724448680877),
    @(0.0429132841527462, 1.122159481048584, 0.1867720484733581),
    @(0.06475171446800231, -1.842216849327088, -0.6108652353286743),
    @(0.0862847194075584, -5.713422775268555, -1.346194267272949),
    @(-0.1818851232528686, -4.851491928100586, 1.392772793769836),
    @(-0.3181080818176269, -3.869678497314453, 0.9886853694915771),
    @(0.1050688251852989, -2.216677188873291, 0.3729332387447357),
    @(0.1996002197265625, 1.434922456741333, -0.2237294018268585)
)

$row = 23
foreach ($vals in $newRows) {
    $ws.Range("A$row").Value = $vals[0]
    $ws.Range("B$row").Value = $vals[1]
    $ws.Range("C$row").Value = $vals[2]
    $row++
}
